$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add 6 new rows (16-21) to the translation table, continuing the existing
# "Название файла / Номер строки / Англ. Строка / Переведенная строка /
#  Конвертированная строка" layout.
#
# Row formatting follows the same repeating pattern used throughout the
# sheet:
#   - "header" rows with only column A filled in use the style of row 10
#   - rows that start a new source-file block use the style of row 3
#   - continuation rows (shorter, single line) use the style of row 6
# We copy (Paste Special -> formats only) from those template rows so the
# cell styles/borders match exactly, then fill in the values, and finally
# set each row's height the same way Excel does when a row is sized to fit
# its (wrapped) contents.
# ---------------------------------------------------------------------------

# Row 16: file-name-only row (like row 10)
$ws.Range("A10:E10").Copy() | Out-Null
$ws.Range("A16:E16").PasteSpecial(-4122) | Out-Null

# Row 17: new file block, first line (like row 3)
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A17:E17").PasteSpecial(-4122) | Out-Null

# Row 18: new file block, first line (like row 6)
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A18:E18").PasteSpecial(-4122) | Out-Null

# Row 19: continuation row (like row 10)
$ws.Range("A10:E10").Copy() | Out-Null
$ws.Range("A19:E19").PasteSpecial(-4122) | Out-Null

# Row 20: continuation row (like row 3)
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A20:E20").PasteSpecial(-4122) | Out-Null

# Row 21: new file block, first line (like row 6)
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A21:E21").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Fill in the cell values. Values are written in the same left-to-right,
# row-by-row order the translator originally typed them in (this also
# determines the order new entries land in xl/sharedStrings.xml), with the
# B column (plain numbers, not shared strings) poked in alongside.
# ---------------------------------------------------------------------------

# Row 16
$ws.Range("A16").Value = "SCRIPT/T01P02A/um2503.ssb"

# Row 17
$ws.Range("C17").Value = " Putting your heart into a gift of\ntreasure...[K] Sounds nice, huh?"
$ws.Range("A17").Value = "SCRIPT/D73P23A/us3104.ssb"
$ws.Range("B17").Value = 212
$ws.Range("D17").Value = " Вложить сердце и душу в\nподарок-сокровище...[K] Звучит недурно, а?"
$ws.Range("E17").Value = " Âìïçéóû òåñäøå é äôšô â\nðïäàñïë-òïëñïâéþå...[K] Èâôœéó îåäôñîï, à?"

# Row 18
$ws.Range("A18").Value = "SCRIPT/D73P26A/us3107.ssb"
$ws.Range("B18").Value = 196
$ws.Range("C18").Value = " You know... Team [CS:X]Frontier[CR] is \njust awesome."

# Row 19 (interleaved with row 18 the way the translator actually typed it)
$ws.Range("C19").Value = " They\'re carrying heavy goods and\nestablishing a base for everyone."
$ws.Range("B19").Value = 199

$ws.Range("D18").Value = " Знаешь... Команда [CS:X]Рубеж[CR]\nвеликолепна."
$ws.Range("D19").Value = " Они таскают тяжёлые припасы и\nразбивают для всех лагеря."
$ws.Range("E18").Value = " Èîàåšû... Ëïíàîäà [CS:X]Ñôáåç[CR]\nâåìéëïìåðîà."
$ws.Range("E19").Value = " Ïîé óàòëàýó óÿçæìúå ðñéðàòú é\nñàèáéâàýó äìÿ âòåö ìàãåñÿ."

# Row 20
$ws.Range("B20").Value = 183
$ws.Range("C20").Value = " I wonder if Team [CS:X]Frontier[CR] has\nreached the top yet…"
$ws.Range("D20").Value = " Интересно, Команда [CS:X]Рубеж[CR] уже\nдостигла вершины?"
$ws.Range("E20").Value = " Éîóåñåòîï, Ëïíàîäà [CS:X]Ñôáåç[CR] ôçå\näïòóéãìà âåñšéîú?"

# Row 21
$ws.Range("C21").Value = " When you\'re having trouble,\ntaking a break and resting could help."
$ws.Range("A21").Value = "SCRIPT/P01P04A/us3116.ssb"
$ws.Range("B21").Value = 170
$ws.Range("D21").Value = " Если что-то не получается,\nвозможно, стоит сделать перерыв и\nотдохнуть."
$ws.Range("E21").Value = " Åòìé œóï-óï îå ðïìôœàåóòÿ,\nâïèíïçîï, òóïéó òäåìàóû ðåñåñúâ é\nïóäïöîôóû."

# ---------------------------------------------------------------------------
# Row heights: rows 16/17/18/21 hold two-line wrapped text (43.2pt, same as
# the other "full" rows); rows 19/20 hold a single line of wrapped text
# (21.6pt, same as the other "continuation" rows). Row 15 grows from one to
# two lines as well, so its custom height increases to match.
# ---------------------------------------------------------------------------
$ws.Rows.Item(15).RowHeight = 43.2
$ws.Rows.Item(16).RowHeight = 43.2
$ws.Rows.Item(17).RowHeight = 43.2
$ws.Rows.Item(18).RowHeight = 43.2
$ws.Rows.Item(19).RowHeight = 21.6
$ws.Rows.Item(20).RowHeight = 21.6
$ws.Rows.Item(21).RowHeight = 43.2

# ---------------------------------------------------------------------------
# Scroll / selection, matching where the author ended up after the edit
# ---------------------------------------------------------------------------
$ws.Range("E20").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 19
